$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = -0.01064165722506156; "C2" = 1.4457030009673;    "D2" = 8.757661946936523; "E2" = 2.959334713569339; "F2" = 3.028956062860261; "G2" = 22
    "B3" = -0.1057000273176439;  "C3" = 1.429381664828222;  "D3" = 8.348262251371553; "E3" = 2.889335953358756; "F3" = 2.958706515248924; "G3" = 21
    "B4" = -0.5533781653781482;  "C4" = 0.9838110950134815; "D4" = 4.162027514816797; "E4" = 2.04010478035242;  "F4" = 2.014630727793039; "G4" = 20
    "B5" = -0.1254169797412491;  "C5" = 0.6292378930413424; "D5" = 0.94278273826223;  "E5" = 0.9709699986416831;"F5" = 0.9892200370224752;"G5" = 19
    "B6" = -0.07984804025652048; "C6" = 0.6307496203848055; "D6" = 0.7016782114047131;"E6" = 0.837662349282044; "F6" = 0.8580225224554666;"G6" = 18
    "B7" = -0.12596358335141;    "C7" = 0.5576815445666784; "D7" = 0.6625899158906151;"E7" = 0.8139962628234942;"F7" = 0.8289410623092455;"G7" = 17
    "B8" = 0.004716166735411878; "C8" = 0.4246577403850764; "D8" = 0.3914765290368093;"E8" = 0.6256808523814752;"F8" = 0.6461820480808862;"G8" = 16
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
